$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New student rows (3-7) are appended below the existing single record.
#
# Cells are written one-by-one, in a specific order, so that the workbook's
# shared-string table is built up in the same sequence as in the source
# data: first the manually-entered row 3 (car, apell, nom, cum, prom),
# then column-by-column batches for the rows that were bulk-pasted in
# (column A for rows 4-6, the name columns, row 7 as a whole, then the
# "cum" column for rows 4-7).
#
# A handful of values in the "cum"/"prom" columns look like numbers
# ("7.76", "7.5", ...) but must be stored as literal text, matching the
# source workbook. A leading apostrophe forces Excel to keep such values
# as text instead of silently converting them to numbers; the style is
# then reset to "Normal" so the cell itself keeps the default (no "s"
# attribute), matching the target file.

function Set-TextCell($addr, $val) {
    $ws.Range($addr).Value = $val
}

function Set-NumericLookingTextCell($addr, $val) {
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).Style = "Normal"
}

# Row 3 - Bryan Lobos (entered first, manually)
Set-TextCell "A3" "LC13004"
Set-TextCell "C3" "Lobos"
Set-TextCell "B3" "Bryan"
Set-NumericLookingTextCell "F3" "7.76"
Set-NumericLookingTextCell "H3" "7.52"

# Column A for rows 4-6
Set-TextCell "A4" "SR11038"
Set-TextCell "A5" "AM11098"
Set-TextCell "A6" "BV13003"

# Names for rows 4-6 (nom/apell)
Set-TextCell "B4" "Rodrigo"
Set-TextCell "C4" "Segovia"
Set-TextCell "C5" "Motto"
Set-TextCell "B5" "Dario"
Set-TextCell "B6" "Elias"
Set-TextCell "C6" "Barrera"

# Row 7 entered as a whole
Set-TextCell "A7" "LL13002"
Set-TextCell "B7" "Alam"
Set-TextCell "C7" "Lopez"

# "cum" column for rows 4-7
Set-NumericLookingTextCell "F4" "7.5"
Set-NumericLookingTextCell "F5" "7.6"
Set-NumericLookingTextCell "F6" "8.3"
Set-NumericLookingTextCell "F7" "7.5"

# Remaining plain numeric cells
$ws.Range("D3").Value = 32
$ws.Range("E3").Value = 0
$ws.Range("G3").Value = 2013
$ws.Range("I3").Value = 1

$ws.Range("D4").Value = 34
$ws.Range("E4").Value = 2
$ws.Range("G4").Value = 2011
$ws.Range("H4").Value = 8
$ws.Range("I4").Value = 1

$ws.Range("D5").Value = 33
$ws.Range("E5").Value = 2
$ws.Range("G5").Value = 2011
$ws.Range("H5").Value = 7
$ws.Range("I5").Value = 1

$ws.Range("D6").Value = 32
$ws.Range("E6").Value = 0
$ws.Range("G6").Value = 2013
$ws.Range("H6").Value = 9
$ws.Range("I6").Value = 1

$ws.Range("D7").Value = 32
$ws.Range("E7").Value = 0
$ws.Range("G7").Value = 2013
$ws.Range("H7").Value = 7
$ws.Range("I7").Value = 1

# Leave the cursor on H3, matching the workbook's saved selection.
$ws.Range("H3").Select() | Out-Null
